# Weekly driver report update for 2025-04-19
# Refresh the "Good Drivers" table (rows 12-17) on the Driver Summary sheet:
# client counts / percentages were recalculated for the new sampling window
# and the rows were re-ordered to reflect the new ranking. Two driver
# versions (21.60.2.1 and 22.50.1.1) no longer have a resolvable "Driver
# Vintage" date this week, so column E is cleared for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B12").Value = 56018
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = ""

# Row 13: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B13").Value = 34244
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = ""

# Row 14: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 442178
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").Value = "'2024-11-10"

# Row 15: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B15").Value = 77849
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").Value = "'2021-08-18"

# Row 16: Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = "'2020-08-05"

# Row 17: Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "'2020-01-06"

# The literal date strings in E14:E17 pick up a stray quote-prefixed /
# date-flavoured number format when written via COM. Re-stamp each with
# the plain "General" formatting the column already used (borrowed from
# its own D-column neighbour, which keeps style untouched) so the look
# of the column matches the rest of the table.
$ws.Range("D14").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
